# Case_1_238 vm_pu.xlsx update: bus 0 (slack) voltage set-point lowered
# from 1.05 p.u. to 1.02 p.u. (380 kV case) and downstream bus voltages
# recomputed by the load-flow solver for every result row (2-25).
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$row = New-Object 'object[,]' 1,13
$row[0,0] = 1.02
$row[0,1] = 1.036368143126336
$row[0,2] = 1.053840555119002
$row[0,3] = 1.035015862455223
$row[0,4] = 1.059513628349498
$row[0,5] = 1
$row[0,6] = $null
$row[0,7] = 1.040660984946148
$row[0,8] = 1.041476841805171
$row[0,9] = 1.056585413056771
$row[0,10] = 1.037813678117573
$row[0,11] = 1.062242926509055
$row[0,12] = 1.042955858379597
$ws.Range("B2:N2").Value = $row

$row = New-Object 'object[,]' 1,13
$row[0,0] = 1.02
$row[0,1] = 1.039169436755782
$row[0,2] = 1.054866274878997
$row[0,3] = 1.037468291508507
$row[0,4] = 1.061023536434782
$row[0,5] = 1
$row[0,6] = $null
$row[0,7] = 1.040958147004656
$row[0,8] = 1.043912761824634
$row[0,9] = 1.057423713373109
$row[0,10] = 1.040071049511232
$row[0,11] = 1.063565319930882
$row[0,12] = 1.045395237684893
$ws.Range("B3:N3").Value = $row

$row = New-Object 'object[,]' 1,13
$row[0,0] = 1.02
$row[0,1] = 1.04097124072823
$row[0,2] = 1.055525711607635
$row[0,3] = 1.039045235623069
$row[0,4] = 1.061994299336337
$row[0,5] = 1
$row[0,6] = $null
$row[0,7] = 1.041147290137623
$row[0,8] = 1.045478412541795
$row[0,9] = 1.057961280360791
$row[0,10] = 1.041521523547177
$row[0,11] = 1.06441425697275
$row[0,12] = 1.046963111805654
$ws.Range("B4:N4").Value = $row

$row = New-Object 'object[,]' 1,13
$row[0,0] = 1.02
$row[0,1] = 1.041726200812968
$row[0,2] = 1.055801932227526
$row[0,3] = 1.039705865005467
$row[0,4] = 1.062400940888096
$row[0,5] = 1
$row[0,6] = $null
$row[0,7] = 1.041226060926908
$row[0,8] = 1.046134148983072
$row[0,9] = 1.058186123643426
$row[0,10] = 1.042128920803723
$row[0,11] = 1.064769564121847
$row[0,12] = 1.047619779467892
$ws.Range("B5:N5").Value = $row

$row = New-Object 'object[,]' 1,13
$row[0,0] = 1.02
$row[0,1] = 1.041852816332024
$row[0,2] = 1.055848252326231
$row[0,3] = 1.039816653447023
$row[0,4] = 1.062469132584908
$row[0,5] = 1
$row[0,6] = $null
$row[0,7] = 1.041239243400419
$row[0,8] = 1.046244107394361
$row[0,9] = 1.058223808846106
$row[0,10] = 1.04223076748807
$row[0,11] = 1.064829129545226
$row[0,12] = 1.047729894032737
$ws.Range("B6:N6").Value = $row

$row = New-Object 'object[,]' 1,13
$row[0,0] = 1.02
$row[0,1] = 1.040981338341472
$row[0,2] = 1.055529406414446
$row[0,3] = 1.039054072010716
$row[0,4] = 1.061999738627277
$row[0,5] = 1
$row[0,6] = $null
$row[0,7] = 1.04114834559482
$row[0,8] = 1.045487184118843
$row[0,9] = 1.057964289224761
$row[0,10] = 1.04152964890058
$row[0,11] = 1.064419010799338
$row[0,12] = 1.046971895839347
$ws.Range("B7:N7").Value = $row

$row = New-Object 'object[,]' 1,13
$row[0,0] = 1.02
$row[0,1] = 1.037317137677723
$row[0,2] = 1.054188096477892
$row[0,3] = 1.035846768527656
$row[0,4] = 1.060025220157762
$row[0,5] = 1
$row[0,6] = $null
$row[0,7] = 1.040762066723179
$row[0,8] = 1.042302295723492
$row[0,9] = 1.056869738460644
$row[0,10] = 1.038578714388291
$row[0,11] = 1.062691247873938
$row[0,12] = 1.043782484537163
$ws.Range("B8:N8").Value = $row

$row = New-Object 'object[,]' 1,13
$row[0,0] = 1.02
$row[0,1] = 1.030774184558118
$row[0,2] = 1.051791068963442
$row[0,3] = 1.030116115054306
$row[0,4] = 1.056496694810925
$row[0,5] = 1
$row[0,6] = $null
$row[0,7] = 1.040057025486971
$row[0,8] = 1.036606412258363
$row[0,9] = 1.054902980331625
$row[0,10] = 1.033298007049391
$row[0,11] = 1.059593828144082
$row[0,12] = 1.038078512264188
$ws.Range("B9:N9").Value = $row

$row = New-Object 'object[,]' 1,13
$row[0,0] = 1.02
$row[0,1] = 1.026349603748548
$row[0,2] = 1.050169484285686
$row[0,3] = 1.026238550741007
$row[0,4] = 1.054109411706165
$row[0,5] = 1
$row[0,6] = $null
$row[0,7] = 1.039570163344461
$row[0,8] = 1.032748778464308
$row[0,9] = 1.053565209297826
$row[0,10] = 1.02971939224454
$row[0,11] = 1.057491531223697
$row[0,12] = 1.034215400187669
$ws.Range("B10:N10").Value = $row

$row = New-Object 'object[,]' 1,13
$row[0,0] = 1.02
$row[0,1] = 1.024417796027874
$row[0,2] = 1.049461487154069
$row[0,3] = 1.024545055710583
$row[0,4] = 1.05306697912015
$row[0,5] = 1
$row[0,6] = $null
$row[0,7] = 1.039355250500342
$row[0,8] = 1.031063126596132
$row[0,9] = 1.052979393765907
$row[0,10] = 1.028155150669328
$row[0,11] = 1.056571944655592
$row[0,12] = 1.032527354500492
$ws.Range("B11:N11").Value = $row

$row = New-Object 'object[,]' 1,13
$row[0,0] = 1.02
$row[0,1] = 1.023697749423475
$row[0,2] = 1.049197606626415
$row[0,3] = 1.023913759969086
$row[0,4] = 1.052678425194593
$row[0,5] = 1
$row[0,6] = $null
$row[0,7] = 1.039274797461853
$row[0,8] = 1.030434624494366
$row[0,9] = 1.052760790588172
$row[0,10] = 1.027571840296985
$row[0,11] = 1.056228939346492
$row[0,12] = 1.031897959853652
$ws.Range("B12:N12").Value = $row

$row = New-Object 'object[,]' 1,13
$row[0,0] = 1.02
$row[0,1] = 1.023852315939266
$row[0,2] = 1.049254250922083
$row[0,3] = 1.024049278462287
$row[0,4] = 1.05276183294591
$row[0,5] = 1
$row[0,6] = $null
$row[0,7] = 1.039292083336239
$row[0,8] = 1.030569549180882
$row[0,9] = 1.052807727541151
$row[0,10] = 1.027697066837338
$row[0,11] = 1.056302580445017
$row[0,12] = 1.032033076148698
$ws.Range("B13:N13").Value = $row

$row = New-Object 'object[,]' 1,13
$row[0,0] = 1.02
$row[0,1] = 1.024358328037092
$row[0,2] = 1.049439693170197
$row[0,3] = 1.024492919107017
$row[0,4] = 1.05303488885176
$row[0,5] = 1
$row[0,6] = $null
$row[0,7] = 1.039348613027272
$row[0,8] = 1.031011223365237
$row[0,9] = 1.052961344587767
$row[0,10] = 1.028106981040408
$row[0,11] = 1.056543621111257
$row[0,12] = 1.032475377561056
$ws.Range("B14:N14").Value = $row

$row = New-Object 'object[,]' 1,13
$row[0,0] = 1.02
$row[0,1] = 1.02466976638444
$row[0,2] = 1.049553830499458
$row[0,3] = 1.024765959335432
$row[0,4] = 1.053202947839709
$row[0,5] = 1
$row[0,6] = $null
$row[0,7] = 1.039383359761478
$row[0,8] = 1.031283036092519
$row[0,9] = 1.053055859286793
$row[0,10] = 1.028359238087067
$row[0,11] = 1.056691943651057
$row[0,12] = 1.032747576293588
$ws.Range("B15:N15").Value = $row

$row = New-Object 'object[,]' 1,13
$row[0,0] = 1.02
$row[0,1] = 1.026477467836904
$row[0,2] = 1.050216346859611
$row[0,3] = 1.026350630466422
$row[0,4] = 1.054178407558362
$row[0,5] = 1
$row[0,6] = $null
$row[0,7] = 1.039584339283102
$row[0,8] = 1.032860320964347
$row[0,9] = 1.053603948137482
$row[0,10] = 1.029822890100818
$row[0,11] = 1.057552362626884
$row[0,12] = 1.034327101090851
$ws.Range("B16:N16").Value = $row

$row = New-Object 'object[,]' 1,13
$row[0,0] = 1.02
$row[0,1] = 1.027607059944457
$row[0,2] = 1.050630347314755
$row[0,3] = 1.027340719022908
$row[0,4] = 1.054787926057839
$row[0,5] = 1
$row[0,6] = $null
$row[0,7] = 1.039709304835415
$row[0,8] = 1.033845564643092
$row[0,9] = 1.053945981125957
$row[0,10] = 1.030737017491203
$row[0,11] = 1.058089571599384
$row[0,12] = 1.035313743928628
$ws.Range("B17:N17").Value = $row

$row = New-Object 'object[,]' 1,13
$row[0,0] = 1.02
$row[0,1] = 1.028264400820402
$row[0,2] = 1.050871264752575
$row[0,3] = 1.02791682894172
$row[0,4] = 1.055142608800217
$row[0,5] = 1
$row[0,6] = $null
$row[0,7] = 1.039781800327325
$row[0,8] = 1.034418772803411
$row[0,9] = 1.054144852461236
$row[0,10] = 1.031268801588888
$row[0,11] = 1.058402023777175
$row[0,12] = 1.035887766110284
$ws.Range("B18:N18").Value = $row

$row = New-Object 'object[,]' 1,13
$row[0,0] = 1.02
$row[0,1] = 1.028488280159614
$row[0,2] = 1.050953316684825
$row[0,3] = 1.028113033831396
$row[0,4] = 1.055263405247028
$row[0,5] = 1
$row[0,6] = $null
$row[0,7] = 1.039806452723931
$row[0,8] = 1.0346139754504
$row[0,9] = 1.054212556089927
$row[0,10] = 1.031449889184376
$row[0,11] = 1.058508411549404
$row[0,12] = 1.036083245967421
$ws.Range("B19:N19").Value = $row

$row = New-Object 'object[,]' 1,13
$row[0,0] = 1.02
$row[0,1] = 1.02748602443538
$row[0,2] = 1.050585987270004
$row[0,3] = 1.027234636444034
$row[0,4] = 1.054722617549969
$row[0,5] = 1
$row[0,6] = $null
$row[0,7] = 1.039695938115547
$row[0,8] = 1.03374000972676
$row[0,9] = 1.053909349610825
$row[0,10] = 1.030639086759044
$row[0,11] = 1.058032026756989
$row[0,12] = 1.035208039112208
$ws.Range("B20:N20").Value = $row

$row = New-Object 'object[,]' 1,13
$row[0,0] = 1.02
$row[0,1] = 1.024209389505149
$row[0,2] = 1.049385110038344
$row[0,3] = 1.0243623409334
$row[0,4] = 1.052954518211768
$row[0,5] = 1
$row[0,6] = $null
$row[0,7] = 1.039331983765251
$row[0,8] = 1.030881227580175
$row[0,9] = 1.052916136115265
$row[0,10] = 1.027986335124256
$row[0,11] = 1.056472680409018
$row[0,12] = 1.032345197167069
$ws.Range("B21:N21").Value = $row

$row = New-Object 'object[,]' 1,13
$row[0,0] = 1.02
$row[0,1] = 1.022134800981735
$row[0,2] = 1.048624859237018
$row[0,3] = 1.022543319336172
$row[0,4] = 1.051835025307164
$row[0,5] = 1
$row[0,6] = $null
$row[0,7] = 1.039099531397159
$row[0,8] = 1.029070009308005
$row[0,9] = 1.052285836813247
$row[0,10] = 1.026305206184227
$row[0,11] = 1.055483965919672
$row[0,12] = 1.030531406757201
$ws.Range("B22:N22").Value = $row

$row = New-Object 'object[,]' 1,13
$row[0,0] = 1.02
$row[0,1] = 1.023235979914197
$row[0,2] = 1.04902838388331
$row[0,3] = 1.02350888572479
$row[0,4] = 1.052429243570317
$row[0,5] = 1
$row[0,6] = $null
$row[0,7] = 1.039223105021404
$row[0,8] = 1.030031505193941
$row[0,9] = 1.052620529804457
$row[0,10] = 1.027197685295515
$row[0,11] = 1.0560089003842
$row[0,12] = 1.031494268077582
$ws.Range("B23:N23").Value = $row

$row = New-Object 'object[,]' 1,13
$row[0,0] = 1.02
$row[0,1] = 1.027540719897988
$row[0,2] = 1.050606033400535
$row[0,3] = 1.027282574892046
$row[0,4] = 1.054752130244517
$row[0,5] = 1
$row[0,6] = $null
$row[0,7] = 1.039701979180416
$row[0,8] = 1.033787709980979
$row[0,9] = 1.053925903763723
$row[0,10] = 1.030683341797199
$row[0,11] = 1.058058031551068
$row[0,12] = 1.03525580710626
$ws.Range("B24:N24").Value = $row

$row = New-Object 'object[,]' 1,13
$row[0,0] = 1.02
$row[0,1] = 1.032476380426144
$row[0,2] = 1.052414832012875
$row[0,3] = 1.03160739417433
$row[0,4] = 1.057414920177241
$row[0,5] = 1
$row[0,6] = $null
$row[0,7] = 1.040242228989817
$row[0,8] = 1.038089269934865
$row[0,9] = 1.055416045038124
$row[0,10] = 1.034673158184794
$row[0,11] = 1.060401034170698
$row[0,12] = 1.039563475768676
$ws.Range("B25:N25").Value = $row
